# case1/5/Qlm2_2.xlsx : refresh computed values in columns A & B (rows 1-32)
# and slightly widen columns A and B.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1, 1).Value = -0.30188945895206132
$ws.Cells.Item(1, 2).Value = 0.30096135179041994
$ws.Cells.Item(2, 1).Value = -0.17336280217121036
$ws.Cells.Item(2, 2).Value = 0.17177231031311901
$ws.Cells.Item(3, 1).Value = -0.12205368009854922
$ws.Cells.Item(3, 2).Value = 0.12157007454450053
$ws.Cells.Item(4, 1).Value = -0.11357007462281388
$ws.Cells.Item(4, 2).Value = 0.11313637403569388
$ws.Cells.Item(5, 1).Value = -0.11013637408313404
$ws.Cells.Item(5, 2).Value = 0.10866300450737754
$ws.Cells.Item(6, 1).Value = -0.0093780570885204639
$ws.Cells.Item(6, 2).Value = 0.0093163740029726938
$ws.Cells.Item(7, 1).Value = 0.0006836258839375553
$ws.Cells.Item(7, 2).Value = -0.00068805772522395614
$ws.Cells.Item(8, 1).Value = 0.01068805761231495
$ws.Cells.Item(8, 2).Value = -0.010693808263849025
$ws.Cells.Item(9, 1).Value = 0.012693808210830326
$ws.Cells.Item(9, 2).Value = -0.012705212609569916
$ws.Cells.Item(10, 1).Value = 0.014705212558515868
$ws.Cells.Item(10, 2).Value = -0.014705148636787513
$ws.Cells.Item(11, 1).Value = -0.024390409828590087
$ws.Cells.Item(11, 2).Value = 0.024364320594403743
$ws.Cells.Item(12, 1).Value = -0.020864320657016933
$ws.Cells.Item(12, 2).Value = 0.02067001349584574
$ws.Cells.Item(13, 1).Value = -0.017170013562787467
$ws.Cells.Item(13, 2).Value = 0.01708172581470091
$ws.Cells.Item(14, 1).Value = -0.0090817259155135943
$ws.Cells.Item(14, 2).Value = 0.0090530555355403663
$ws.Cells.Item(15, 1).Value = -0.0080530555860258701
$ws.Cells.Item(15, 2).Value = 0.0080345587504879745
$ws.Cells.Item(16, 1).Value = -0.0060345588089387725
$ws.Cells.Item(16, 2).Value = 0.0060032980499524768
$ws.Cells.Item(17, 1).Value = -0.004003298109339859
$ws.Cells.Item(17, 2).Value = 0.0039999999257585017
$ws.Cells.Item(18, 1).Value = -0.067245767723431271
$ws.Cells.Item(18, 2).Value = 0.067060876476677578
$ws.Cells.Item(19, 1).Value = -0.012091669813892647
$ws.Cells.Item(19, 2).Value = 0.012016349921941671
$ws.Cells.Item(20, 1).Value = -0.0080163499573799868
$ws.Cells.Item(20, 2).Value = 0.0080055709555022503
$ws.Cells.Item(21, 1).Value = -0.0040055709913398019
$ws.Cells.Item(21, 2).Value = 0.0039999999639039885
$ws.Cells.Item(22, 1).Value = -0.045718630257990611
$ws.Cells.Item(22, 2).Value = 0.045503515755328294
$ws.Cells.Item(23, 1).Value = -0.040503515808726576
$ws.Cells.Item(23, 2).Value = 0.040099840363637895
$ws.Cells.Item(24, 1).Value = -0.020099840530540369
$ws.Cells.Item(24, 2).Value = 0.019999999830710991
$ws.Cells.Item(25, 1).Value = -0.097284947473767147
$ws.Cells.Item(25, 2).Value = 0.097159046042419206
$ws.Cells.Item(26, 1).Value = -0.094659046102949063
$ws.Cells.Item(26, 2).Value = 0.094496369497225174
$ws.Cells.Item(27, 1).Value = -0.091996369561179847
$ws.Cells.Item(27, 2).Value = 0.091029462761073709
$ws.Cells.Item(28, 1).Value = -0.089029462837141971
$ws.Cells.Item(28, 2).Value = 0.088367550746105472
$ws.Cells.Item(29, 1).Value = -0.08136755086967451
$ws.Cells.Item(29, 2).Value = 0.08117588811132137
$ws.Cells.Item(30, 1).Value = -0.021175888626699724
$ws.Cells.Item(30, 2).Value = 0.021023318686021142
$ws.Cells.Item(31, 1).Value = -0.014023318819070596
$ws.Cells.Item(31, 2).Value = 0.014000919784969312
$ws.Cells.Item(32, 1).Value = -0.0040009199404185125
$ws.Cells.Item(32, 2).Value = 0.0039999998880997367

$ws.Range("A:A").ColumnWidth = 14.75
$ws.Range("B:B").ColumnWidth = 15.583333333333336

